$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matching original inlineStr type) for price cells
# whose new value would otherwise be auto-parsed as a number.
$textCells = @("D5", "D6", "D14", "D19", "D21", "D22", "D24", "D28", "D29", "D30", "D31", "D32", "D33", "D36", "D37", "D40", "D41", "D43", "D44", "D46", "D47")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range('D2').Value = '63.199.08'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '2.614.83'
$ws.Range('E3').Value = '  -2.38%  '
$ws.Range('D5').Value = '606.38'
$ws.Range('E5').Value = '  +1.01%  '
$ws.Range('D6').Value = '145.62'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D9').Value = '2.613.66'
$ws.Range('E9').Value = '  -2.24%  '
$ws.Range('E10').Value = '  +0.29%  '
$ws.Range('E11').Value = '  -3.44%  '
$ws.Range('E12').Value = '  +4.26%  '
$ws.Range('E13').Value = '  -0.55%  '
$ws.Range('D14').Value = '27.18'
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('D15').Value = '3.082.59'
$ws.Range('E15').Value = '  -2.42%  '
$ws.Range('D16').Value = '63.064.95'
$ws.Range('E16').Value = '  -0.84%  '
$ws.Range('D18').Value = '2.629.17'
$ws.Range('E18').Value = '  -2.03%  '
$ws.Range('D19').Value = '11.52'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('D21').Value = '342.05'
$ws.Range('E21').Value = '  +0.56%  '
$ws.Range('D22').Value = '6.87'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D24').Value = '5.74'
$ws.Range('E24').Value = '  -1.30%  '
$ws.Range('E25').Value = '  -2.62%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +2.92%  '
$ws.Range('D28').Value = '9.02'
$ws.Range('E28').Value = '  +5.42%  '
$ws.Range('D29').Value = '546.30'
$ws.Range('E29').Value = '  +1.00%  '
$ws.Range('D30').Value = '0.161'
$ws.Range('E30').Value = '  -2.58%  '
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('D32').Value = '7.82'
$ws.Range('E32').Value = '  -0.92%  '
$ws.Range('D33').Value = '2.02'
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('E34').Value = '  +2.51%  '
$ws.Range('E35').Value = '  -5.76%  '
$ws.Range('D36').Value = '5.24'
$ws.Range('E36').Value = '  +1.49%  '
$ws.Range('D37').Value = '168.77'
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('E39').Value = '  -1.11%  '
$ws.Range('D40').Value = '1.94'
$ws.Range('E40').Value = '  +5.13%  '
$ws.Range('D41').Value = '18.90'
$ws.Range('E41').Value = '  -1.85%  '
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').Value = '165.78'
$ws.Range('E43').Value = '  -4.88%  '
$ws.Range('D44').Value = '39.65'
$ws.Range('E44').Value = '  -1.39%  '
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('D46').Value = '21.73'
$ws.Range('E46').Value = '  -2.69%  '
$ws.Range('D47').Value = '0.0562'
$ws.Range('E47').Value = '  -0.75%  '
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('E50').Value = '  -1.10%  '
$ws.Range('E51').Value = '  +10.26%  '
